$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the "Then_Goto" / "Else_Goto" headers to "Then_Question" / "Else_Question"
# in the Survey data dictionary (Sheet1 header row).
$thenCell = $ws.Cells.Find("Then_Goto")
if ($thenCell -ne $null) {
    $thenCell.Value = "Then_Question"
}

$elseCell = $ws.Cells.Find("Else_Goto")
if ($elseCell -ne $null) {
    $elseCell.Value = "Else_Question"
}

# Reflect the author's resulting selection on Sheet1.
[void]$ws.Range("L10").Select()
